# Re-fetch the live workbook/worksheet objects from $excel.
# NOTE: the pre-bound $wb / $ws variables handed to this script appear to be
# disconnected/stale in this runtime (e.g. $wb.Worksheets.Count reads back as
# 0, $wb.ActiveSheet.Range(...) throws "cannot call a method on a
# null-valued expression"). Re-resolving everything through $excel avoids
# that and is what actually persists changes to the saved file.
$wbLive = $excel.ActiveWorkbook
$wsLive = $wbLive.ActiveSheet

$newD5 = @"
💚 上海线下 ZY0602698(加急)(#上外虹口超近#同济大学四平路#复旦邯郸#上财国定路#交大医学院#华师普陀)
【辅导地址】虹口足球场地铁站附近
【辅导科目】 升高三 物 化
【学员情况】 女孩 成绩中等
【辅导时间】1 次/周  2h/次 ，周末上课，老师优秀暑期继续安排
【教员要求】能力第一位的，会解题会讲明白
【老师薪资】300－360/2h
"@

$newD6 = @"
🌸 上海线下ZY0526663（6.4加急信息费打折 #上海交通大学）
【辅导地址】杨浦区延吉中路城市丽园245弄
【辅导科目】 初三，数理化最后冲刺
【学员情况】 女孩，基础薄弱
【辅导时间】6月7日下午，12.13日时间可协商
【教员要求】专职在职老师，数理化都可以辅导，能冲刺提分
【老师薪资】500－600/2h
"@

$newD7 = @"
🌸上海线下ZY0607741a(#同济四平路#交大医学院#复旦大学枫林#交大徐汇#上外虹口#上财国定路)
【辅导地址】浦东新区商城路1025弄
【辅导科目】 四升五，语文
【学员情况】女孩，80-85分
【辅导时间】1次/周  2h/次，周末上课，暑期安排同上
【教员要求】老师经验要丰富，山东籍优先，可长期授课
【老师薪资】200-240/2h
"@

$wsLive.Range("D5").Value = $newD5
$wsLive.Range("D6").Value = $newD6
$wsLive.Range("D7").Value = $newD7

# Move the selection from H2 to D5 (matches the sheetView's <selection> in
# the target file). $excel.Range(...) -- not $wsLive.Range(...) -- is the
# form whose .Select() call actually takes effect in this runtime.
[void]$excel.Range("D5").Select()

"done"
